# Weekly update: insert a new week's Jengibre price row at row 21,
# pushing the existing rows 21-84 down to 22-85, then populate the
# new row 21 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 21 (shifts 21..84 -> 22..85)
$ws.Rows.Item(21).Insert()

# Fill the new row 21 with the latest week's record
$ws.Cells.Item(21, 1).Value = 8
$ws.Cells.Item(21, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(21, 3).Value = "Coquimbo"
$ws.Cells.Item(21, 4).Value = 44914
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 100114007
$ws.Cells.Item(21, 7).Value = "Jengibre"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 460
$ws.Cells.Item(21, 11).Value = 14000
$ws.Cells.Item(21, 12).Value = 15000
$ws.Cells.Item(21, 13).Value = 14500
$ws.Cells.Item(21, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(21, 15).Value = "Perú"
$ws.Cells.Item(21, 16).Value = 1115
$ws.Cells.Item(21, 17).Value = 13
$ws.Cells.Item(21, 18).Value = "Hortaliza"
